$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 207 (the COCOOH reaction row) - remaining rows shift up automatically
$ws.Rows.Item(207).Delete()

# Update sheet view: selection and top-left cell as recorded after the edit
$ws.Range("A2:A227").Select()
$excel.ActiveWindow.ScrollRow = 211
